$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.Execute($find, $true, $false, $false, $false, $false,
                         $true, 1, $false, $replace, 2)
}

# Title changes in heading and bold "Play Lock it Link Night Life for Free - Game Review"
# appears twice with identical text -> replace all occurrences
Replace-Text "Play Lock it Link Night Life for Free - Game Review" "Play Lock it Link Night Life Free - Review"

# "What we like" bullet list
Replace-Text "High payout percentage of 96.2%" "High payout percentage"
Replace-Text "Progressive jackpots that can be won by players" "Progressive jackpots"
Replace-Text "Well-designed symbols related to wealth and luxury" "Well-designed symbols"
Replace-Text "Mini-games that provide free spins" "Similar game available"

# "What we don't like" bullet list
Replace-Text "Limited theme variety" "Limited number of mini-games"
Replace-Text "Limited variety in terms of bonus features" "No bonus rounds"

# Meta description (italic)
Replace-Text "Read our comprehensive review on Lock it Link Night Life slot game. Play it now for free and win big with high payout percentage and progressive jackpots." "Read our review of Lock it Link Night Life to experience the thrill of playing for free."
